$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 67-79 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A67:A79").NumberFormat = "@"
$ws.Range("A67").Value = '2026-01-28'
$ws.Range("B67").Value = '16:15:58'
$ws.Range("C67").Value = '16:00'
$ws.Range("D67").Value = 'Bathroom'
$ws.Range("E67").Value = 'No Motion'
$ws.Range("F67").Value = 'Inactive'
$ws.Range("A68").Value = '2026-01-28'
$ws.Range("B68").Value = '16:16:00'
$ws.Range("C68").Value = '16:00'
$ws.Range("D68").Value = 'Bathroom'
$ws.Range("E68").Value = 'No Motion'
$ws.Range("F68").Value = 'Inactive'
$ws.Range("A69").Value = '2026-01-28'
$ws.Range("B69").Value = '16:16:04'
$ws.Range("C69").Value = '16:00'
$ws.Range("D69").Value = 'Bathroom'
$ws.Range("E69").Value = 'No Motion'
$ws.Range("F69").Value = 'Inactive'
$ws.Range("A70").Value = '2026-01-28'
$ws.Range("B70").Value = '16:16:09'
$ws.Range("C70").Value = '16:00'
$ws.Range("D70").Value = 'Bathroom'
$ws.Range("E70").Value = 'No Motion'
$ws.Range("F70").Value = 'Inactive'
$ws.Range("A71").Value = '2026-01-28'
$ws.Range("B71").Value = '16:16:14'
$ws.Range("C71").Value = '16:00'
$ws.Range("D71").Value = 'Bathroom'
$ws.Range("E71").Value = 'No Motion'
$ws.Range("F71").Value = 'Inactive'
$ws.Range("A72").Value = '2026-01-28'
$ws.Range("B72").Value = '16:16:19'
$ws.Range("C72").Value = '16:00'
$ws.Range("D72").Value = 'Bathroom'
$ws.Range("E72").Value = 'No Motion'
$ws.Range("F72").Value = 'Inactive'
$ws.Range("A73").Value = '2026-01-28'
$ws.Range("B73").Value = '16:16:24'
$ws.Range("C73").Value = '16:00'
$ws.Range("D73").Value = 'Bathroom'
$ws.Range("E73").Value = 'No Motion'
$ws.Range("F73").Value = 'Inactive'
$ws.Range("A74").Value = '2026-01-28'
$ws.Range("B74").Value = '16:16:29'
$ws.Range("C74").Value = '16:00'
$ws.Range("D74").Value = 'Bathroom'
$ws.Range("E74").Value = 'No Motion'
$ws.Range("F74").Value = 'Inactive'
$ws.Range("A75").Value = '2026-01-28'
$ws.Range("B75").Value = '16:16:34'
$ws.Range("C75").Value = '16:00'
$ws.Range("D75").Value = 'Bathroom'
$ws.Range("E75").Value = 'No Motion'
$ws.Range("F75").Value = 'Inactive'
$ws.Range("A76").Value = '2026-01-28'
$ws.Range("B76").Value = '16:16:39'
$ws.Range("C76").Value = '16:00'
$ws.Range("D76").Value = 'Bathroom'
$ws.Range("E76").Value = 'No Motion'
$ws.Range("F76").Value = 'Inactive'
$ws.Range("A77").Value = '2026-01-28'
$ws.Range("B77").Value = '16:16:44'
$ws.Range("C77").Value = '16:00'
$ws.Range("D77").Value = 'Bathroom'
$ws.Range("E77").Value = 'No Motion'
$ws.Range("F77").Value = 'Inactive'
$ws.Range("A78").Value = '2026-01-28'
$ws.Range("B78").Value = '16:16:49'
$ws.Range("C78").Value = '16:00'
$ws.Range("D78").Value = 'Bathroom'
$ws.Range("E78").Value = 'No Motion'
$ws.Range("F78").Value = 'Inactive'
$ws.Range("A79").Value = '2026-01-28'
$ws.Range("B79").Value = '16:16:55'
$ws.Range("C79").Value = '16:00'
$ws.Range("D79").Value = 'Bathroom'
$ws.Range("E79").Value = 'No Motion'
$ws.Range("F79").Value = 'Inactive'

# --- Humidity sheet: append rows 66-78 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A66:A78").NumberFormat = "@"
$ws.Range("E66:E78").NumberFormat = "@"
$ws.Range("A66").Value = '2026-01-28'
$ws.Range("B66").Value = '16:15:59'
$ws.Range("C66").Value = '16:00'
$ws.Range("D66").Value = 'Bathroom'
$ws.Range("E66").Value = '88.3%'
$ws.Range("F66").Value = 'Active'
$ws.Range("A67").Value = '2026-01-28'
$ws.Range("B67").Value = '16:15:59'
$ws.Range("C67").Value = '16:00'
$ws.Range("D67").Value = 'Bathroom'
$ws.Range("E67").Value = '88.3%'
$ws.Range("F67").Value = 'Active'
$ws.Range("A68").Value = '2026-01-28'
$ws.Range("B68").Value = '16:16:07'
$ws.Range("C68").Value = '16:00'
$ws.Range("D68").Value = 'Bathroom'
$ws.Range("E68").Value = '88.3%'
$ws.Range("F68").Value = 'Active'
$ws.Range("A69").Value = '2026-01-28'
$ws.Range("B69").Value = '16:16:15'
$ws.Range("C69").Value = '16:00'
$ws.Range("D69").Value = 'Bathroom'
$ws.Range("E69").Value = '87.4%'
$ws.Range("F69").Value = 'Active'
$ws.Range("A70").Value = '2026-01-28'
$ws.Range("B70").Value = '16:16:19'
$ws.Range("C70").Value = '16:00'
$ws.Range("D70").Value = 'Bathroom'
$ws.Range("E70").Value = '88.3%'
$ws.Range("F70").Value = 'Active'
$ws.Range("A71").Value = '2026-01-28'
$ws.Range("B71").Value = '16:16:23'
$ws.Range("C71").Value = '16:00'
$ws.Range("D71").Value = 'Bathroom'
$ws.Range("E71").Value = '87.4%'
$ws.Range("F71").Value = 'Active'
$ws.Range("A72").Value = '2026-01-28'
$ws.Range("B72").Value = '16:16:27'
$ws.Range("C72").Value = '16:00'
$ws.Range("D72").Value = 'Bathroom'
$ws.Range("E72").Value = '88.3%'
$ws.Range("F72").Value = 'Active'
$ws.Range("A73").Value = '2026-01-28'
$ws.Range("B73").Value = '16:16:35'
$ws.Range("C73").Value = '16:00'
$ws.Range("D73").Value = 'Bathroom'
$ws.Range("E73").Value = '88.3%'
$ws.Range("F73").Value = 'Active'
$ws.Range("A74").Value = '2026-01-28'
$ws.Range("B74").Value = '16:16:39'
$ws.Range("C74").Value = '16:00'
$ws.Range("D74").Value = 'Bathroom'
$ws.Range("E74").Value = '88.3%'
$ws.Range("F74").Value = 'Active'
$ws.Range("A75").Value = '2026-01-28'
$ws.Range("B75").Value = '16:16:43'
$ws.Range("C75").Value = '16:00'
$ws.Range("D75").Value = 'Bathroom'
$ws.Range("E75").Value = '87.4%'
$ws.Range("F75").Value = 'Active'
$ws.Range("A76").Value = '2026-01-28'
$ws.Range("B76").Value = '16:16:47'
$ws.Range("C76").Value = '16:00'
$ws.Range("D76").Value = 'Bathroom'
$ws.Range("E76").Value = '86.9%'
$ws.Range("F76").Value = 'Active'
$ws.Range("A77").Value = '2026-01-28'
$ws.Range("B77").Value = '16:16:51'
$ws.Range("C77").Value = '16:00'
$ws.Range("D77").Value = 'Bathroom'
$ws.Range("E77").Value = '88.3%'
$ws.Range("F77").Value = 'Active'
$ws.Range("A78").Value = '2026-01-28'
$ws.Range("B78").Value = '16:16:55'
$ws.Range("C78").Value = '16:00'
$ws.Range("D78").Value = 'Bathroom'
$ws.Range("E78").Value = '87.4%'
$ws.Range("F78").Value = 'Active'

# --- Temperature sheet: append rows 66-78 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A66:A78").NumberFormat = "@"
$ws.Range("A66").Value = '2026-01-28'
$ws.Range("B66").Value = '16:15:59'
$ws.Range("C66").Value = '16:00'
$ws.Range("D66").Value = 'Bathroom'
$ws.Range("E66").Value = '22.8C'
$ws.Range("F66").Value = 'Active'
$ws.Range("A67").Value = '2026-01-28'
$ws.Range("B67").Value = '16:15:59'
$ws.Range("C67").Value = '16:00'
$ws.Range("D67").Value = 'Bathroom'
$ws.Range("E67").Value = '22.8C'
$ws.Range("F67").Value = 'Active'
$ws.Range("A68").Value = '2026-01-28'
$ws.Range("B68").Value = '16:16:07'
$ws.Range("C68").Value = '16:00'
$ws.Range("D68").Value = 'Bathroom'
$ws.Range("E68").Value = '22.8C'
$ws.Range("F68").Value = 'Active'
$ws.Range("A69").Value = '2026-01-28'
$ws.Range("B69").Value = '16:16:15'
$ws.Range("C69").Value = '16:00'
$ws.Range("D69").Value = 'Bathroom'
$ws.Range("E69").Value = '22.8C'
$ws.Range("F69").Value = 'Active'
$ws.Range("A70").Value = '2026-01-28'
$ws.Range("B70").Value = '16:16:19'
$ws.Range("C70").Value = '16:00'
$ws.Range("D70").Value = 'Bathroom'
$ws.Range("E70").Value = '22.8C'
$ws.Range("F70").Value = 'Active'
$ws.Range("A71").Value = '2026-01-28'
$ws.Range("B71").Value = '16:16:23'
$ws.Range("C71").Value = '16:00'
$ws.Range("D71").Value = 'Bathroom'
$ws.Range("E71").Value = '22.8C'
$ws.Range("F71").Value = 'Active'
$ws.Range("A72").Value = '2026-01-28'
$ws.Range("B72").Value = '16:16:27'
$ws.Range("C72").Value = '16:00'
$ws.Range("D72").Value = 'Bathroom'
$ws.Range("E72").Value = '22.8C'
$ws.Range("F72").Value = 'Active'
$ws.Range("A73").Value = '2026-01-28'
$ws.Range("B73").Value = '16:16:35'
$ws.Range("C73").Value = '16:00'
$ws.Range("D73").Value = 'Bathroom'
$ws.Range("E73").Value = '22.7C'
$ws.Range("F73").Value = 'Active'
$ws.Range("A74").Value = '2026-01-28'
$ws.Range("B74").Value = '16:16:39'
$ws.Range("C74").Value = '16:00'
$ws.Range("D74").Value = 'Bathroom'
$ws.Range("E74").Value = '22.7C'
$ws.Range("F74").Value = 'Active'
$ws.Range("A75").Value = '2026-01-28'
$ws.Range("B75").Value = '16:16:43'
$ws.Range("C75").Value = '16:00'
$ws.Range("D75").Value = 'Bathroom'
$ws.Range("E75").Value = '22.7C'
$ws.Range("F75").Value = 'Active'
$ws.Range("A76").Value = '2026-01-28'
$ws.Range("B76").Value = '16:16:47'
$ws.Range("C76").Value = '16:00'
$ws.Range("D76").Value = 'Bathroom'
$ws.Range("E76").Value = '22.8C'
$ws.Range("F76").Value = 'Active'
$ws.Range("A77").Value = '2026-01-28'
$ws.Range("B77").Value = '16:16:51'
$ws.Range("C77").Value = '16:00'
$ws.Range("D77").Value = 'Bathroom'
$ws.Range("E77").Value = '22.8C'
$ws.Range("F77").Value = 'Active'
$ws.Range("A78").Value = '2026-01-28'
$ws.Range("B78").Value = '16:16:55'
$ws.Range("C78").Value = '16:00'
$ws.Range("D78").Value = 'Bathroom'
$ws.Range("E78").Value = '22.7C'
$ws.Range("F78").Value = 'Active'
